$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '37.260.08'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '2.066.96'
$ws.Range('E3').Value = '  -0.48%  '
$ws.Range('E4').Value = '  +0.00%  '
Set-TextValue 'D5' '233.55'
$ws.Range('E5').Value = '  -1.29%  '
Set-TextValue 'D6' '0.623'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('E7').Value = '  +0.04%  '
Set-TextValue 'D8' '56.83'
$ws.Range('E8').Value = '  -2.52%  '
Set-TextValue 'D9' '0.382'
$ws.Range('E9').Value = '  -0.77%  '
Set-TextValue 'D10' '0.0763'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').Value = '2.371.22'
$ws.Range('E12').Value = '  -0.53%  '
Set-TextValue 'D13' '14.63'
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('E15').Value = '  -0.62%  '
Set-TextValue 'D16' '5.13'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').Value = '2.069.52'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '37.231.32'
$ws.Range('E18').Value = '  -1.52%  '
Set-TextValue 'D19' '6.38'
$ws.Range('E19').Value = '  +3.57%  '
Set-TextValue 'D20' '69.44'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('E21').Value = '  -0.28%  '
Set-TextValue 'D22' '226.64'
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('E23').Value = '  +0.16%  '
Set-TextValue 'D24' '2.43'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  -2.65%  '
Set-TextValue 'D26' '166.70'
$ws.Range('E26').Value = '  +1.89%  '
Set-TextValue 'D27' '8.77'
$ws.Range('E27').Value = '  -0.95%  '
Set-TextValue 'D28' '1.43'
$ws.Range('E28').Value = '  +2.73%  '
Set-TextValue 'D29' '19.07'
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('E30').Value = '  -4.31%  '
Set-TextValue 'D31' '0.117'
$ws.Range('E31').Value = '  -1.10%  '
Set-TextValue 'D32' '4.46'
$ws.Range('E32').Value = '  -0.59%  '
Set-TextValue 'D33' '0.0615'
$ws.Range('E33').Value = '  -2.59%  '
$ws.Range('E34').Value = '  +1.83%  '
$ws.Range('E35').Value = '  -4.57%  '
$ws.Range('E36').Value = '  -0.02%  '
Set-TextValue 'D37' '1.76'
$ws.Range('E37').Value = '  -1.12%  '
Set-TextValue 'D38' '3.23'
$ws.Range('E38').Value = '  -3.68%  '
$ws.Range('E39').Value = '  -4.42%  '
Set-TextValue 'D40' '2.95'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('E41').Value = '  +3.76%  '
$ws.Range('D42').Value = '1.480.58'
$ws.Range('E42').Value = '  -0.22%  '
Set-TextValue 'D43' '96.06'
$ws.Range('E43').Value = '  +0.64%  '
Set-TextValue 'D44' '0.0935'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('E45').Value = '  +3.20%  '
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('E47').Value = '  -0.86%  '
Set-TextValue 'D48' '15.06'
$ws.Range('E48').Value = '  -8.54%  '
Set-TextValue 'D49' '7.15'
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('E50').Value = '  +1.00%  '
$ws.Range('D51').Value = '2.260.27'
$ws.Range('E51').Value = '  -0.46%  '
